# Applies the cryptos-list refresh described by the commit:
#   "Updated cryptos list on Wed Mar 29 06:32:15 UTC 2023 with GitHub Actions"
#
# Column D ("Price") values are forced to text (leading single-quote = Excel's
# literal-text marker) because many of them look like numbers (e.g. "0.9970",
# "1.130", "5.710") and a plain numeric assignment would silently lose trailing
# zeros / exact digits through floating point conversion - the source workbook
# stores these as text. The per-cell Style reset below strips the "quote prefix"
# flag Excel attaches when literal text is entered, so the cells end up with no
# style override at all (matching the original workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ref -> new value. Kind "T" marks values that must be forced to text.
$updates = @(
    @{ Ref = 'D2'; Value = '27.880.73'; Kind = 'T' }
    @{ Ref = 'E2'; Value = '  +2.98%  '; Kind = 'S' }
    @{ Ref = 'D3'; Value = '1.802.67'; Kind = 'T' }
    @{ Ref = 'E3'; Value = '  +4.09%  '; Kind = 'S' }
    @{ Ref = 'D4'; Value = '0.9945'; Kind = 'T' }
    @{ Ref = 'E4'; Value = '  -0.54%  '; Kind = 'S' }
    @{ Ref = 'D5'; Value = '314.81'; Kind = 'T' }
    @{ Ref = 'E5'; Value = '  +1.37%  '; Kind = 'S' }
    @{ Ref = 'D6'; Value = '0.9970'; Kind = 'T' }
    @{ Ref = 'E6'; Value = '  -0.21%  '; Kind = 'S' }
    @{ Ref = 'D7'; Value = '0.5446'; Kind = 'T' }
    @{ Ref = 'E7'; Value = '  +11.83%  '; Kind = 'S' }
    @{ Ref = 'D8'; Value = '0.3788'; Kind = 'T' }
    @{ Ref = 'E8'; Value = '  +8.03%  '; Kind = 'S' }
    @{ Ref = 'B9'; Value = 'Dogecoin'; Kind = 'S' }
    @{ Ref = 'C9'; Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; Kind = 'S' }
    @{ Ref = 'D9'; Value = '0.07555'; Kind = 'T' }
    @{ Ref = 'E9'; Value = '  +3.68%  '; Kind = 'S' }
    @{ Ref = 'B10'; Value = 'OKB'; Kind = 'S' }
    @{ Ref = 'C10'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; Kind = 'S' }
    @{ Ref = 'D10'; Value = '42.65'; Kind = 'T' }
    @{ Ref = 'E10'; Value = '  +1.35%  '; Kind = 'S' }
    @{ Ref = 'D11'; Value = '1.130'; Kind = 'T' }
    @{ Ref = 'E11'; Value = '  +7.16%  '; Kind = 'S' }
    @{ Ref = 'D12'; Value = '0.9989'; Kind = 'T' }
    @{ Ref = 'E12'; Value = '  -0.08%  '; Kind = 'S' }
    @{ Ref = 'D13'; Value = '21.06'; Kind = 'T' }
    @{ Ref = 'E13'; Value = '  +5.25%  '; Kind = 'S' }
    @{ Ref = 'D14'; Value = '6.222'; Kind = 'T' }
    @{ Ref = 'E14'; Value = '  +5.49%  '; Kind = 'S' }
    @{ Ref = 'D15'; Value = '1.788.22'; Kind = 'T' }
    @{ Ref = 'E15'; Value = '  +2.86%  '; Kind = 'S' }
    @{ Ref = 'D16'; Value = '7.121'; Kind = 'T' }
    @{ Ref = 'E16'; Value = '  +3.26%  '; Kind = 'S' }
    @{ Ref = 'D17'; Value = '91.09'; Kind = 'T' }
    @{ Ref = 'E17'; Value = '  +4.32%  '; Kind = 'S' }
    @{ Ref = 'D18'; Value = '0.00001075'; Kind = 'T' }
    @{ Ref = 'E18'; Value = '  +3.30%  '; Kind = 'S' }
    @{ Ref = 'D19'; Value = '0.06500'; Kind = 'T' }
    @{ Ref = 'E19'; Value = '  +1.40%  '; Kind = 'S' }
    @{ Ref = 'D20'; Value = '0.9978'; Kind = 'T' }
    @{ Ref = 'E20'; Value = '  -0.12%  '; Kind = 'S' }
    @{ Ref = 'D21'; Value = '17.09'; Kind = 'T' }
    @{ Ref = 'E21'; Value = '  +3.18%  '; Kind = 'S' }
    @{ Ref = 'D22'; Value = '5.954'; Kind = 'T' }
    @{ Ref = 'E22'; Value = '  +4.72%  '; Kind = 'S' }
    @{ Ref = 'D23'; Value = '27.840.07'; Kind = 'T' }
    @{ Ref = 'E23'; Value = '  +2.59%  '; Kind = 'S' }
    @{ Ref = 'D24'; Value = '11.25'; Kind = 'T' }
    @{ Ref = 'E24'; Value = '  +3.66%  '; Kind = 'S' }
    @{ Ref = 'D25'; Value = '2.103'; Kind = 'T' }
    @{ Ref = 'E25'; Value = '  +1.15%  '; Kind = 'S' }
    @{ Ref = 'D26'; Value = '20.60'; Kind = 'T' }
    @{ Ref = 'E26'; Value = '  +3.03%  '; Kind = 'S' }
    @{ Ref = 'D27'; Value = '155.09'; Kind = 'T' }
    @{ Ref = 'E27'; Value = '  +1.26%  '; Kind = 'S' }
    @{ Ref = 'D28'; Value = '2.398'; Kind = 'T' }
    @{ Ref = 'E28'; Value = '  +14.48%  '; Kind = 'S' }
    @{ Ref = 'D29'; Value = '1.996.56'; Kind = 'T' }
    @{ Ref = 'E29'; Value = '  +3.67%  '; Kind = 'S' }
    @{ Ref = 'D30'; Value = '122.68'; Kind = 'T' }
    @{ Ref = 'E30'; Value = '  +0.59%  '; Kind = 'S' }
    @{ Ref = 'D31'; Value = '1.150'; Kind = 'T' }
    @{ Ref = 'E31'; Value = '  +10.37%  '; Kind = 'S' }
    @{ Ref = 'D32'; Value = '0.1036'; Kind = 'T' }
    @{ Ref = 'E32'; Value = '  +11.11%  '; Kind = 'S' }
    @{ Ref = 'D33'; Value = '5.710'; Kind = 'T' }
    @{ Ref = 'E33'; Value = '  +5.40%  '; Kind = 'S' }
    @{ Ref = 'D34'; Value = '3.617'; Kind = 'T' }
    @{ Ref = 'E34'; Value = '  +0.09%  '; Kind = 'S' }
    @{ Ref = 'D35'; Value = '0.02297'; Kind = 'T' }
    @{ Ref = 'E35'; Value = '  +4.41%  '; Kind = 'S' }
    @{ Ref = 'D36'; Value = '8.661'; Kind = 'T' }
    @{ Ref = 'E36'; Value = '  +15.28%  '; Kind = 'S' }
    @{ Ref = 'B37'; Value = 'InternetComputer(DFINITY)'; Kind = 'S' }
    @{ Ref = 'C37'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Kind = 'S' }
    @{ Ref = 'D37'; Value = '5.037'; Kind = 'T' }
    @{ Ref = 'E37'; Value = '  +5.46%  '; Kind = 'S' }
    @{ Ref = 'B38'; Value = 'Algorand'; Kind = 'S' }
    @{ Ref = 'C38'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Kind = 'S' }
    @{ Ref = 'D38'; Value = '0.2104'; Kind = 'T' }
    @{ Ref = 'E38'; Value = '  +4.92%  '; Kind = 'S' }
    @{ Ref = 'B39'; Value = 'Hedera'; Kind = 'S' }
    @{ Ref = 'C39'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; Kind = 'S' }
    @{ Ref = 'D39'; Value = '0.06052'; Kind = 'T' }
    @{ Ref = 'E39'; Value = '  +2.17%  '; Kind = 'S' }
    @{ Ref = 'D40'; Value = '11.45'; Kind = 'T' }
    @{ Ref = 'E40'; Value = '  +3.72%  '; Kind = 'S' }
    @{ Ref = 'D41'; Value = '0.6292'; Kind = 'T' }
    @{ Ref = 'E41'; Value = '  +4.54%  '; Kind = 'S' }
    @{ Ref = 'D42'; Value = '0.9989'; Kind = 'T' }
    @{ Ref = 'E42'; Value = '  +0.03%  '; Kind = 'S' }
    @{ Ref = 'D43'; Value = '1.398'; Kind = 'T' }
    @{ Ref = 'E43'; Value = '  -2.11%  '; Kind = 'S' }
    @{ Ref = 'D44'; Value = '1.150'; Kind = 'T' }
    @{ Ref = 'E44'; Value = '  +4.98%  '; Kind = 'S' }
    @{ Ref = 'D45'; Value = '13.36'; Kind = 'T' }
    @{ Ref = 'E45'; Value = '  +4.45%  '; Kind = 'S' }
    @{ Ref = 'D46'; Value = '0.5914'; Kind = 'T' }
    @{ Ref = 'E46'; Value = '  +4.15%  '; Kind = 'S' }
    @{ Ref = 'D47'; Value = '3.637'; Kind = 'T' }
    @{ Ref = 'E47'; Value = '  +1.40%  '; Kind = 'S' }
    @{ Ref = 'D48'; Value = '121.88'; Kind = 'T' }
    @{ Ref = 'E48'; Value = '  +2.68%  '; Kind = 'S' }
    @{ Ref = 'D49'; Value = '1.927'; Kind = 'T' }
    @{ Ref = 'E49'; Value = '  +4.26%  '; Kind = 'S' }
    @{ Ref = 'D50'; Value = '1.135'; Kind = 'T' }
    @{ Ref = 'E50'; Value = '  +2.23%  '; Kind = 'S' }
    @{ Ref = 'D51'; Value = '0.06777'; Kind = 'T' }
    @{ Ref = 'E51'; Value = '  +1.75%  '; Kind = 'S' }
)

$quoteChar = [string][char]39   # single-quote / apostrophe

foreach ($u in $updates) {
    $range = $ws.Range($u.Ref)
    if ($u.Kind -eq "T") {
        # Leading apostrophe forces Excel to keep the numeric-looking string as text.
        $range.Value = $quoteChar + $u.Value
    } else {
        $range.Value = $u.Value
    }
}

# Strip the quote-prefix style flag picked up above so the cells carry no style
# override, exactly as in the source file.
foreach ($u in $updates) {
    if ($u.Kind -eq "T") {
        $ws.Range($u.Ref).Style = "Normal"
    }
}
